# The data rows (2-14) of the "Artfynd" sheet were re-sorted upstream.
# Concretely: the records that used to be the last two rows (old row 12 and
# old row 14) now sort to the top (new rows 2 and 3), and every other record
# keeps its relative order, shifting down by two rows.
#
# old row -> new row
#   2->4   3->5   4->6   5->7   6->8   7->9   8->10  9->11
#   10->12 11->13 12->2  13->14 14->3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns I (Antal), Y (Startdatum) and AA (Slutdatum) hold plain text that
# looks like a number / ISO date ("25", "2012-11-09", ...). Excel's COM
# layer auto-converts such text to a real number / date serial on write, so
# force those columns to Text format up front to keep them as literal
# strings, matching the source workbook (Z/AB "HH:MM" time text is left
# alone on write by this engine, so no workaround is needed there).
$ws.Range("I2:I14").NumberFormat = "@"
$ws.Range("Y2:Y14").NumberFormat = "@"
$ws.Range("AA2:AA14").NumberFormat = "@"

# Snapshot every used row (2-14) across all columns that carry data
# (A:AY) BEFORE writing anything back, so reads never see an
# already-overwritten row.
$snapshot = @{}
for ($r = 2; $r -le 14; $r++) {
    $snapshot[$r] = $ws.Range("A$r`:AY$r").Value2
}

$rowMap = @{
    2  = 4
    3  = 5
    4  = 6
    5  = 7
    6  = 8
    7  = 9
    8  = 10
    9  = 11
    10 = 12
    11 = 13
    12 = 2
    13 = 14
    14 = 3
}

foreach ($oldRow in $rowMap.Keys) {
    $newRow = $rowMap[$oldRow]
    $data = $snapshot[$oldRow]
    $dest = $ws.Range("A$newRow`:AY$newRow")
    $dest.Value2 = $data
}
